$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 45863
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 33.703
$ws.Cells.Item(3, 1).Value = 45863.01041666666
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 33.569
$ws.Cells.Item(4, 1).Value = 45863.02083333334
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 28.484
$ws.Cells.Item(5, 1).Value = 45863.03125
$ws.Cells.Item(5, 2).Value = 0.029
$ws.Cells.Item(5, 3).Value = 12.937
$ws.Cells.Item(6, 1).Value = 45863.04166666666
$ws.Cells.Item(6, 2).Value = 1.062
$ws.Cells.Item(6, 3).Value = 0.101
$ws.Cells.Item(7, 1).Value = 45863.05208333334
$ws.Cells.Item(7, 2).Value = 2.014
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = 45863.0625
$ws.Cells.Item(8, 2).Value = 3.477
$ws.Cells.Item(8, 3).Value = 0.089
$ws.Cells.Item(9, 1).Value = 45863.07291666666
$ws.Cells.Item(9, 2).Value = 6.375
$ws.Cells.Item(9, 3).Value = 1.43
$ws.Cells.Item(10, 1).Value = 45863.08333333334
$ws.Cells.Item(10, 2).Value = 1.184
$ws.Cells.Item(10, 3).Value = 1.549
$ws.Cells.Item(11, 1).Value = 45863.09375
$ws.Cells.Item(11, 2).Value = 0.386
$ws.Cells.Item(11, 3).Value = 0.398
$ws.Cells.Item(12, 1).Value = 45863.10416666666
$ws.Cells.Item(12, 2).Value = 1.014
$ws.Cells.Item(12, 3).Value = 0.532
$ws.Cells.Item(13, 1).Value = 45863.11458333334
$ws.Cells.Item(13, 2).Value = 6.831
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = 45863.125
$ws.Cells.Item(14, 2).Value = 2.664
$ws.Cells.Item(14, 3).Value = 1.799
$ws.Cells.Item(15, 1).Value = 45863.13541666666
$ws.Cells.Item(15, 2).Value = 1.836
$ws.Cells.Item(15, 3).Value = 0.391
$ws.Cells.Item(16, 1).Value = 45863.14583333334
$ws.Cells.Item(16, 2).Value = 2.357
$ws.Cells.Item(16, 3).Value = 1.75
$ws.Cells.Item(17, 1).Value = 45863.15625
$ws.Cells.Item(17, 2).Value = 0.002
$ws.Cells.Item(17, 3).Value = 8.218999999999999
$ws.Cells.Item(18, 1).Value = 45863.16666666666
$ws.Cells.Item(18, 2).Value = 0.051
$ws.Cells.Item(18, 3).Value = 2.927
$ws.Cells.Item(19, 1).Value = 45863.17708333334
$ws.Cells.Item(19, 2).Value = 11.876
$ws.Cells.Item(19, 3).Value = 0.055
$ws.Cells.Item(20, 1).Value = 45863.1875
$ws.Cells.Item(20, 2).Value = 9.153
$ws.Cells.Item(20, 3).Value = 0.005
$ws.Cells.Item(21, 1).Value = 45863.19791666666
$ws.Cells.Item(21, 2).Value = 15.631
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(22, 1).Value = 45863.20833333334
$ws.Cells.Item(22, 2).Value = 0.212
$ws.Cells.Item(22, 3).Value = 9.175000000000001
$ws.Cells.Item(23, 1).Value = 45863.21875
$ws.Cells.Item(23, 2).Value = 0.163
$ws.Cells.Item(23, 3).Value = 3.131
$ws.Cells.Item(24, 1).Value = 45863.22916666666
$ws.Cells.Item(24, 2).Value = 1.395
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(25, 1).Value = 45863.23958333334
$ws.Cells.Item(25, 2).Value = 1.359
$ws.Cells.Item(25, 3).Value = 0.024
$ws.Cells.Item(26, 1).Value = 45863.25
$ws.Cells.Item(26, 2).Value = 1.984
$ws.Cells.Item(26, 3).Value = 5.14
$ws.Cells.Item(27, 1).Value = 45863.26041666666
$ws.Cells.Item(27, 2).Value = 3.184
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(28, 1).Value = 45863.27083333334
$ws.Cells.Item(28, 2).Value = 1.375
$ws.Cells.Item(28, 3).Value = 0.047
$ws.Cells.Item(29, 1).Value = 45863.28125
$ws.Cells.Item(29, 2).Value = 0.436
$ws.Cells.Item(29, 3).Value = 1.518
$ws.Cells.Item(30, 1).Value = 45863.29166666666
$ws.Cells.Item(30, 2).Value = 4.417
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = 45863.30208333334
$ws.Cells.Item(31, 2).Value = 3.786
$ws.Cells.Item(31, 3).Value = 0.059
$ws.Cells.Item(32, 1).Value = 45863.3125
$ws.Cells.Item(32, 2).Value = 0.025
$ws.Cells.Item(32, 3).Value = 6.827
$ws.Cells.Item(33, 1).Value = 45863.32291666666
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 24.248
$ws.Cells.Item(34, 1).Value = 45863.33333333334
$ws.Cells.Item(34, 2).Value = 68.098
$ws.Cells.Item(34, 3).Value = 0.096
$ws.Cells.Item(35, 1).Value = 45863.34375
$ws.Cells.Item(35, 2).Value = 41.793
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = 45863.35416666666
$ws.Cells.Item(36, 2).Value = 24.138
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 1).Value = 45863.36458333334
$ws.Cells.Item(37, 2).Value = 4.84
$ws.Cells.Item(37, 3).Value = 1.409
$ws.Cells.Item(38, 1).Value = 45863.375
$ws.Cells.Item(38, 2).Value = 22.69
$ws.Cells.Item(38, 3).Value = 0.344
$ws.Cells.Item(39, 1).Value = 45863.38541666666
$ws.Cells.Item(39, 2).Value = 33.201
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = 45863.39583333334
$ws.Cells.Item(40, 2).Value = 40.641
$ws.Cells.Item(40, 3).Value = 0
